$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtractedData")

$ws.Range("A2").Value = 'Honda Civic Prim proprietar / carte service și istoric verificabil la Honda'
$ws.Range("B2").Value = 'https://www.autovit.ro/anunt/honda-civic-1-5-ID7GWuaF.html#xtor=SEC-81'
$ws.Range("C2").Value = '20 000 €'
$ws.Range("D2").Value = 'Cumpana'
$ws.Range("E2").Value = 'Azi 20:37'

$ws.Range("A3").Value = 'Honda civic model 2'
$ws.Range("B3").Value = 'https://www.olx.ro/d/oferta/honda-civic-model-2-IDfFwpn.html#4e6472c1f5'
$ws.Range("C3").Value = '4 900 €'
$ws.Range("D3").Value = 'Targoviste'
$ws.Range("E3").Value = 'Azi 20:14'

$ws.Range("A4").Value = 'Honda Civic EJ9 1997'
$ws.Range("B4").Value = 'https://www.olx.ro/d/oferta/honda-civic-ej9-1997-IDfA6RK.html#4e6472c1f5'
$ws.Range("C4").Value = '1 700 €'
$ws.Range("D4").Value = 'Satu Mare'
$ws.Range("E4").Value = 'Azi 19:23'

$ws.Range("A5").Value = 'Honda Civic'
$ws.Range("B5").Value = 'https://www.autovit.ro/anunt/honda-civic-ID7GwjVG.html#xtor=SEC-81'
$ws.Range("C5").Value = '8 050 €'
$ws.Range("D5").Value = 'Timisoara'
$ws.Range("E5").Value = 'Azi 19:00'

$ws.Range("A6").Value = 'Honda civic 2010 full, volan dreapta 2.2 diesel'
$ws.Range("B6").Value = 'https://www.olx.ro/d/oferta/honda-civic-2010-full-volan-dreapta-2-2-diesel-IDfFu0Q.html#4e6472c1f5'
$ws.Range("C6").Value = '1 900 €'
$ws.Range("D6").Value = 'Onesti'
$ws.Range("E6").Value = 'Azi 17:42'

$ws.Range("A7").Value = 'Vand Honda Civic'
$ws.Range("B7").Value = 'https://www.olx.ro/d/oferta/vand-honda-civic-IDfBia5.html#4e6472c1f5;promoted'
$ws.Range("C7").Value = '3 650 €'
$ws.Range("D7").Value = 'Timisoara'
$ws.Range("E7").Value = 'Azi 16:51'

$ws.Range("A8").Value = 'Honda Civic honda civic Civic 1.5 Sport Plus LED#PANO#SPUR#KEYLESS#SHZ'
$ws.Range("B8").Value = 'https://www.autovit.ro/anunt/honda-civic-ID7GWu5Z.html#xtor=SEC-81'
$ws.Range("C8").Value = '26 990,39 €'
$ws.Range("D8").Value = 'Otopeni'
$ws.Range("E8").Value = 'Azi 16:34'

$ws.Range("A9").Value = 'Honda Civic honda civic Civic Limousine 1.5 Executive Automatik *Leder*'
$ws.Range("B9").Value = 'https://www.autovit.ro/anunt/honda-civic-ID7GWu5U.html#xtor=SEC-81'
$ws.Range("C9").Value = '18 990,02 €'
$ws.Range("D9").Value = 'Otopeni'
$ws.Range("E9").Value = 'Azi 16:34'

$ws.Range("A10").Value = 'Honda Civic Facelift 5DR 1.8 I-VTEC 140CP Euro 5 134k km'
$ws.Range("B10").Value = 'https://www.olx.ro/d/oferta/honda-civic-facelift-5dr-1-8-i-vtec-140cp-euro-5-134k-km-IDfBhGG.html#4e6472c1f5;promoted'
$ws.Range("C10").Value = '6 250 €'
$ws.Range("D10").Value = 'Cluj-Napoca'
$ws.Range("E10").Value = 'Azi 16:13'

$ws.Range("A11").Value = 'Honda Civic ej2 America coupe inmatriculată'
$ws.Range("B11").Value = 'https://www.olx.ro/d/oferta/honda-civic-ej2-america-coupe-inmatriculata-IDfbh5u.html#4e6472c1f5;promoted'
$ws.Range("C11").Value = '3 250 €'
$ws.Range("D11").Value = 'Pantelimon'
$ws.Range("E11").Value = 'Azi 15:44'
